$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 31, shifting existing rows 31-54 down to 32-55.
$ws.Rows.Item(31).Insert()

# Fill the new row 31 with data (same as former row 31 for columns A-L except D,
# and new values for M-T).
$ws.Range("A31").Value = 10
$ws.Range("B31").Value = "Vega Modelo de Temuco"
$ws.Range("C31").Value = "La Araucanía"
$ws.Range("D31").Value = 44566
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103003
$ws.Range("J31").Value = "Damasco"
$ws.Range("K31").Value = "Modesto"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 55
$ws.Range("N31").Value = 15000
$ws.Range("O31").Value = 15000
$ws.Range("P31").Value = 15000
$ws.Range("Q31").Value = "$/bandeja 10 kilos"
$ws.Range("R31").Value = "Provincia de Quillota"
$ws.Range("S31").Value = 1500
$ws.Range("T31").Value = 10
